$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# D4: Status changes from "Not Started" to "In Progress"
$ws1.Range("D4").Value = "In Progress"

# G4: % Complete changes from 0 to 0.1 (10%)
$ws1.Range("G4").Value = 0.1

# E5: was text "Asap", becomes an actual date value (matches existing date format on the cell)
$ws1.Range("E5").Value = 43003

# I5 / I6 / I13: new notes
$ws1.Range("I5").Value = "Meeting planed in Lausanne"
$ws1.Range("I6").Value = "See with Olivier"
$ws1.Range("I13").Value = "see with Tara"

# Make sheet1 the active sheet/tab and select I13
$ws1.Activate()
$ws1.Range("I13").Select()
